# Add three new parameter rows to the "parameter_values" sheet, mirroring
# the existing pattern used for rows 85-87 (parameter name in col A,
# value in col B, "DUMMY" marker in col C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameter_values")

# Row 87 (squeeze_factor_threshold_sba_did_not_run) gains the same "DUMMY"
# marker in column C that the other rows in this block already have.
$ws.Range("C87").Value = "DUMMY"

$ws.Range("A88").Value = "sensitivity_of_assessment_of_obstructed_labour"
$ws.Range("B88").Value = 0.5
$ws.Range("C88").Value = "DUMMY"

$ws.Range("A89").Value = "sensitivity_of_assessment_of_sepsis"
$ws.Range("B89").Value = 0.5
$ws.Range("C89").Value = "DUMMY"

$ws.Range("A90").Value = "sensitivity_of_assessment_of_uterine_rupture"
$ws.Range("B90").Value = 0.7
$ws.Range("C90").Value = "DUMMY"

# Match formatting of column A used elsewhere in the sheet (vertical-center
# style, same as rows above, e.g. A85:A87).
$ws.Range("A88:A90").VerticalAlignment = -4108

# Update selection / active cell to match the post-edit workbook state
# (entire row 90 selected, as Excel does after typing into the last row
# of a block and pressing the down-arrow / row-select).
$ws.Activate()
$ws.Range("A90:XFD90").Select()
